$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the emoji "statut" markers in column A with plain-text / simpler
# symbol equivalents (the book emojis were rendering incorrectly in Excel):
#   📕 (red book)    -> -3
#   📗 (green book)  -> ✅ (check mark)
#   📙 (orange book) -> +3
#   📘 (blue book)   -> ⚠️ (warning sign)

$a2 = $ws.Cells.Item(2, 1)
$a2.Value = "'-3"
$a2.Style = "Normal"

$ws.Cells.Item(3, 1).Value = "✅"

$a4 = $ws.Cells.Item(4, 1)
$a4.Value = "'+3"
$a4.Style = "Normal"

$ws.Cells.Item(5, 1).Value = "⚠️"
$ws.Cells.Item(6, 1).Value = "⚠️"
$ws.Cells.Item(7, 1).Value = "⚠️"
